$wb = $excel.ActiveWorkbook

# --- Login sheet ---
$wsLogin = $wb.Worksheets.Item("Login")

# B2 / B3 currently carry a redundant "apply fill" style that is identical in
# appearance to the plain centered style already used elsewhere (e.g. A2/A3).
# Nudging the Interior property makes the host re-resolve the cell format and
# collapse it onto that existing style record instead of the redundant one.
$wsLogin.Range("B2").Interior.Pattern = $wsLogin.Range("A2").Interior.Pattern
$wsLogin.Range("B3").Interior.Pattern = $wsLogin.Range("A3").Interior.Pattern

# B4 / C4: replace the sample credentials with the real ones, and turn the
# username into a mailto hyperlink (this also applies the built-in
# "Hyperlink" cell style to B4, as Excel does automatically).
$wsLogin.Range("B4").Value = "ex121935@edpr.com"
$wsLogin.Hyperlinks.Add($wsLogin.Range("B4"), "mailto:ex121935@edpr.com", "", "", "ex121935@edpr.com") | Out-Null
$wsLogin.Range("C4").Value = "EDPAuth201..18"

# The workbook was left with "Login" as the active sheet/tab and C4 selected.
$wsLogin.Activate() | Out-Null
$wsLogin.Range("C4").Select() | Out-Null
